$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-3: round Ost/Nord (Q/R) coordinates to nearest integer
$ws.Range("Q2").Value = 702819
$ws.Range("R2").Value = 7300018
$ws.Range("Q3").Value = 702522
$ws.Range("R3").Value = 7300048

# Row 4
$ws.Range("A4").Value = 111865919
$ws.Range("B4").Value = 95538
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 221941
$ws.Range("F4").Value = "Plattlummer"
$ws.Range("G4").Value = "Lycopodium complanatum"
$ws.Range("H4").Value = "L."
$ws.Range("Q4").Value = 702755
$ws.Range("R4").Value = 7299754

# Row 5
$ws.Range("A5").Value = 111865981
$ws.Range("B5").Value = 90652
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 3100
$ws.Range("F5").Value = "Talltaggsvamp"
$ws.Range("G5").Value = "Bankera fuligineoalba"
$ws.Range("H5").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q5").Value = 702696
$ws.Range("R5").Value = 7299770

# Row 6
$ws.Range("A6").Value = 111866048
$ws.Range("B6").Value = 90682
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 2059
$ws.Range("F6").Value = "Skrovlig taggsvamp"
$ws.Range("G6").Value = "Hydnellum scabrosum"
$ws.Range("H6").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q6").Value = 702750
$ws.Range("R6").Value = 7299800

# Row 7
$ws.Range("A7").Value = 111866265
$ws.Range("B7").Value = 78107
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 6453
$ws.Range("F7").Value = "Vedskivlav"
$ws.Range("G7").Value = "Hertelidea botryosa"
$ws.Range("H7").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q7").Value = 702681
$ws.Range("R7").Value = 7299925

# Row 8
$ws.Range("A8").Value = 111866031
$ws.Range("B8").Value = 78107
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 6453
$ws.Range("F8").Value = "Vedskivlav"
$ws.Range("G8").Value = "Hertelidea botryosa"
$ws.Range("H8").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q8").Value = 702750
$ws.Range("R8").Value = 7299800

# Row 9
$ws.Range("A9").Value = 111865866
$ws.Range("B9").Value = 90652
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 3100
$ws.Range("F9").Value = "Talltaggsvamp"
$ws.Range("G9").Value = "Bankera fuligineoalba"
$ws.Range("H9").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q9").Value = 702753
$ws.Range("R9").Value = 7299802

# Row 10
$ws.Range("A10").Value = 111865263
$ws.Range("B10").Value = 90658
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 4361
$ws.Range("F10").Value = "Orange taggsvamp"
$ws.Range("G10").Value = "Hydnellum aurantiacum"
$ws.Range("H10").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q10").Value = 702714
$ws.Range("R10").Value = 7299724

# Row 11
$ws.Range("A11").Value = 111866065
$ws.Range("B11").Value = 78107
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6453
$ws.Range("F11").Value = "Vedskivlav"
$ws.Range("G11").Value = "Hertelidea botryosa"
$ws.Range("H11").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q11").Value = 702768
$ws.Range("R11").Value = 7299828

# Row 12
$ws.Range("A12").Value = 111865488
$ws.Range("B12").Value = 90660
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 4362
$ws.Range("F12").Value = "Blå taggsvamp"
$ws.Range("G12").Value = "Hydnellum caeruleum"
$ws.Range("H12").Value = "(Hornem.) P.Karst."
$ws.Range("Q12").Value = 702716
$ws.Range("R12").Value = 7299725

# Row 13
$ws.Range("A13").Value = 111866170
$ws.Range("B13").Value = 90682
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 2059
$ws.Range("F13").Value = "Skrovlig taggsvamp"
$ws.Range("G13").Value = "Hydnellum scabrosum"
$ws.Range("H13").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q13").Value = 702754
$ws.Range("R13").Value = 7299887

# Row 14
$ws.Range("A14").Value = 111865524
$ws.Range("B14").Value = 90660
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 4362
$ws.Range("F14").Value = "Blå taggsvamp"
$ws.Range("G14").Value = "Hydnellum caeruleum"
$ws.Range("H14").Value = "(Hornem.) P.Karst."
$ws.Range("Q14").Value = 702731
$ws.Range("R14").Value = 7299742

# Row 15
$ws.Range("A15").Value = 111866131
$ws.Range("B15").Value = 90682
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 2059
$ws.Range("F15").Value = "Skrovlig taggsvamp"
$ws.Range("G15").Value = "Hydnellum scabrosum"
$ws.Range("H15").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q15").Value = 702757
$ws.Range("R15").Value = 7299855

# Row 16
$ws.Range("A16").Value = 111866159
$ws.Range("B16").Value = 90652
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 3100
$ws.Range("F16").Value = "Talltaggsvamp"
$ws.Range("G16").Value = "Bankera fuligineoalba"
$ws.Range("H16").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q16").Value = 702755
$ws.Range("R16").Value = 7299865

# Row 17
$ws.Range("A17").Value = 111865578
$ws.Range("B17").Value = 90854
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 2079
$ws.Range("F17").Value = "Nordtagging"
$ws.Range("G17").Value = "Odonticium romellii"
$ws.Range("H17").Value = "(S.Lundell) Parmasto"
$ws.Range("Q17").Value = 702742
$ws.Range("R17").Value = 7299746

# Row 18
$ws.Range("A18").Value = 111866021
$ws.Range("B18").Value = 78107
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 6453
$ws.Range("F18").Value = "Vedskivlav"
$ws.Range("G18").Value = "Hertelidea botryosa"
$ws.Range("H18").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q18").Value = 702738
$ws.Range("R18").Value = 7299806

# Row 19
$ws.Range("A19").Value = 111865668
$ws.Range("B19").Value = 78107
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6453
$ws.Range("F19").Value = "Vedskivlav"
$ws.Range("G19").Value = "Hertelidea botryosa"
$ws.Range("H19").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q19").Value = 702741
$ws.Range("R19").Value = 7299744

# Row 20
$ws.Range("A20").Value = 111866276
$ws.Range("B20").Value = 78107
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6453
$ws.Range("F20").Value = "Vedskivlav"
$ws.Range("G20").Value = "Hertelidea botryosa"
$ws.Range("H20").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q20").Value = 702661
$ws.Range("R20").Value = 7299929

# Row 21
$ws.Range("A21").Value = 111866194
$ws.Range("B21").Value = 90682
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 2059
$ws.Range("F21").Value = "Skrovlig taggsvamp"
$ws.Range("G21").Value = "Hydnellum scabrosum"
$ws.Range("H21").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q21").Value = 702687
$ws.Range("R21").Value = 7299920

# Row 22
$ws.Range("A22").Value = 111865961
$ws.Range("B22").Value = 77267
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6446
$ws.Range("F22").Value = "Kolflarnlav"
$ws.Range("G22").Value = "Carbonicola anthracophila"
$ws.Range("H22").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q22").Value = 702714
$ws.Range("R22").Value = 7299790

# Public-comment cell (AC) moves with the record it belongs to:
# the witch-ring note moves from rows 9/22 to rows 21/15.
$ws.Range("AC9").ClearContents()
$ws.Range("AC22").ClearContents()
$ws.Range("AC15").Value = "Flera fruktkoppar som växer i en häxring"
$ws.Range("AC21").Value = "Flera fruktkoppar som växer i en häxring"
